$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 2223.5386
$ws.Range("J112").Value = 2223.5386
$ws.Range("L112").Value = 6670.6158
$ws.Range("N112").Value = -8886.6158
# Row 114
$ws.Range("H114").Value = 60000
$ws.Range("J114").Value = 60000
$ws.Range("L114").Value = 60000
$ws.Range("N114").Value = -68678
# Row 135
$ws.Range("H135").Value = 761
$ws.Range("I135").Value = 291.4375
$ws.Range("J135").Value = 4517.5
$ws.Range("K135").Value = 2622.9375
$ws.Range("L135").Value = 40657.5
$ws.Range("M135").Value = -87.9375
$ws.Range("N135").Value = -45727.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 555.0645
$ws.Range("I2").Value = 486.44827
$ws.Range("J2").Value = 1550
$ws.Range("K2").Value = 486.44827
$ws.Range("L2").Value = 1550
$ws.Range("M2").Value = -373.44827
$ws.Range("N2").Value = -1776
# Row 32
$ws.Range("H32").Value = 4249.66
$ws.Range("I32").Value = 4328.0205
$ws.Range("J32").Value = 410
$ws.Range("K32").Value = 4328.0205
$ws.Range("L32").Value = 410
$ws.Range("M32").Value = -4041.0205
$ws.Range("N32").Value = -984
# Row 116
$ws.Range("H116").Value = 555.0645
$ws.Range("I116").Value = 486.44827
$ws.Range("J116").Value = 1550
$ws.Range("K116").Value = 486.44827
$ws.Range("L116").Value = 1550
$ws.Range("M116").Value = 1807.55173
$ws.Range("N116").Value = -6138

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 555.0645
$ws.Range("I3").Value = 486.44827
$ws.Range("J3").Value = 1550
$ws.Range("K3").Value = 486.44827
$ws.Range("L3").Value = 1550
$ws.Range("M3").Value = -372.44827
$ws.Range("N3").Value = -1778
# Row 20
$ws.Range("H20").Value = 26326546
$ws.Range("I20").Value = 33345490
$ws.Range("J20").Value = 5499.5
$ws.Range("K20").Value = 33345490
$ws.Range("L20").Value = 5499.5
$ws.Range("M20").Value = -33345243
$ws.Range("N20").Value = -5993.5
# Row 80
$ws.Range("H80").Value = 278
$ws.Range("J80").Value = 350.375
$ws.Range("L80").Value = 350.375
$ws.Range("N80").Value = -2346.375
# Row 83
$ws.Range("H83").Value = 278
$ws.Range("J83").Value = 350.375
$ws.Range("L83").Value = 1751.875
$ws.Range("N83").Value = -11735.875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1704.275
$ws.Range("I16").Value = 1636
$ws.Range("K16").Value = 1636
$ws.Range("M16").Value = -1349
# Row 31
$ws.Range("H31").Value = 4380.885
$ws.Range("I31").Value = 4270.2
$ws.Range("J31").Value = 4531.8184
$ws.Range("K31").Value = 4270.2
$ws.Range("L31").Value = 4531.8184
$ws.Range("M31").Value = -3975.2
$ws.Range("N31").Value = -5121.8184
# Row 34
$ws.Range("H34").Value = 4380.885
$ws.Range("I34").Value = 4270.2
$ws.Range("J34").Value = 4531.8184
$ws.Range("K34").Value = 4270.2
$ws.Range("L34").Value = 4531.8184
$ws.Range("M34").Value = -4068.2
$ws.Range("N34").Value = -4935.8184
# Row 59
$ws.Range("H59").Value = 88862.75
$ws.Range("I59").Value = 80000
$ws.Range("J59").Value = 91817
$ws.Range("K59").Value = 80000
$ws.Range("L59").Value = 91817
$ws.Range("M59").Value = -78855
$ws.Range("N59").Value = -94107
# Row 94
$ws.Range("H94").Value = 2023
$ws.Range("J94").Value = 2474.6365
$ws.Range("L94").Value = 2474.6365
$ws.Range("N94").Value = -3376.6365
# Row 113
$ws.Range("H113").Value = 1704.275
$ws.Range("I113").Value = 1636
$ws.Range("K113").Value = 1636
$ws.Range("M113").Value = 534
# Row 117
$ws.Range("H117").Value = 59797.5
$ws.Range("J117").Value = 59797.5
$ws.Range("L117").Value = 59797.5
$ws.Range("N117").Value = -68975.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1867.8182
$ws.Range("I5").Value = 527
$ws.Range("J5").Value = 2262.1765
$ws.Range("K5").Value = 1581
$ws.Range("L5").Value = 6786.529500000001
$ws.Range("M5").Value = -1469
$ws.Range("N5").Value = -7010.529500000001
# Row 131
$ws.Range("H131").Value = 2015.7826
$ws.Range("I131").Value = 1773.8572
$ws.Range("J131").Value = 2392.111
$ws.Range("K131").Value = 5321.571599999999
$ws.Range("L131").Value = 7176.333
$ws.Range("M131").Value = -281.5715999999993
$ws.Range("N131").Value = -17256.333
# Row 135
$ws.Range("H135").Value = 1867.8182
$ws.Range("I135").Value = 527
$ws.Range("J135").Value = 2262.1765
$ws.Range("K135").Value = 4743
$ws.Range("L135").Value = 20359.5885
$ws.Range("M135").Value = -2208
$ws.Range("N135").Value = -25429.5885

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 8007
$ws.Range("I9").Value = 8007
$ws.Range("K9").Value = 8007
$ws.Range("M9").Value = -7837
# Row 70
$ws.Range("H70").Value = 84559.24000000001
$ws.Range("I70").Value = 147374.22
$ws.Range("K70").Value = 147374.22
$ws.Range("M70").Value = -147104.22
# Row 73
$ws.Range("H73").Value = 84559.24000000001
$ws.Range("I73").Value = 147374.22
$ws.Range("K73").Value = 147374.22
$ws.Range("M73").Value = -146438.22
# Row 132
$ws.Range("H132").Value = 1968.9445
$ws.Range("I132").Value = 1478.3334
$ws.Range("J132").Value = 2950.1667
$ws.Range("K132").Value = 4435.0002
$ws.Range("L132").Value = 8850.500100000001
$ws.Range("M132").Value = -1905.0002
$ws.Range("N132").Value = -13910.5001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3895.7
$ws.Range("J7").Value = 5808.5713
$ws.Range("L7").Value = 5808.5713
$ws.Range("N7").Value = -6032.5713
# Row 17
$ws.Range("H17").Value = 15000
$ws.Range("J17").Value = 15000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15340
# Row 36
$ws.Range("H36").Value = 100000
$ws.Range("J36").Value = 100000
$ws.Range("L36").Value = 100000
$ws.Range("N36").Value = -101124
# Row 68
$ws.Range("H68").Value = 5499
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
# Row 71
$ws.Range("H71").Value = 5499
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
# Row 100
$ws.Range("H100").Value = 3930
$ws.Range("I100").Value = 3925.4707
$ws.Range("K100").Value = 3925.4707
$ws.Range("M100").Value = -3384.4707
# Row 126
$ws.Range("H126").Value = 3895.7
$ws.Range("J126").Value = 5808.5713
$ws.Range("L126").Value = 17425.7139
$ws.Range("N126").Value = -22365.7139
# Row 132
$ws.Range("H132").Value = 3352.7834
$ws.Range("I132").Value = 2630.804
$ws.Range("J132").Value = 7444
$ws.Range("K132").Value = 7892.412
$ws.Range("L132").Value = 22332
$ws.Range("M132").Value = -5362.412
$ws.Range("N132").Value = -27392

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 20
$ws.Range("H20").Value = 14983
$ws.Range("I20").Value = 14950
$ws.Range("J20").Value = 14999.5
$ws.Range("K20").Value = 14950
$ws.Range("L20").Value = 14999.5
$ws.Range("M20").Value = -14710
$ws.Range("N20").Value = -15479.5
# Row 70
$ws.Range("H70").Value = 34872.25
$ws.Range("I70").Value = 10000
$ws.Range("J70").Value = 43163
$ws.Range("K70").Value = 10000
$ws.Range("L70").Value = 43163
$ws.Range("M70").Value = -9685
$ws.Range("N70").Value = -43793
# Row 73
$ws.Range("H73").Value = 34872.25
$ws.Range("I73").Value = 10000
$ws.Range("J73").Value = 43163
$ws.Range("K73").Value = 10000
$ws.Range("L73").Value = 43163
$ws.Range("M73").Value = -8908
$ws.Range("N73").Value = -45347
# Row 109
$ws.Range("H109").Value = 100000
$ws.Range("J109").Value = 100000
$ws.Range("L109").Value = 100000
$ws.Range("N109").Value = -102774

# ---- Deletions (cells removed entirely) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()
